# Add a new "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, value + formatting (reuse the same header style as the
# other columns, e.g. G1, by copying its format onto H1).
$ws.Cells.Item(1, 8).Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data cell for row 2.
$ws.Cells.Item(2, 8).Value = 1
